# Add the new claim/incident row (row 46) to the "NEW" sheet, matching the
# automated export format used by the rest of the table: columns A-L are
# plain text (even when the text looks numeric, e.g. case ids, OTs, dates
# kept as literal strings) and columns M/N (coordinates) are real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Force the text columns to Text format first so numeric-looking strings
# (case id, comuna, OT, attachments flag) and the date string aren't
# auto-converted into numbers/dates by Excel's type inference.
$ws.Range("A$row`:L$row").NumberFormat = "@"

$ws.Range("A$row").Value = "-478"
$ws.Range("B$row").Value = "6/15/2025"
$ws.Range("C$row").Value = "Chivilcoy 4875"
$ws.Range("D$row").Value = "11"
$ws.Range("E$row").Value = "807508509"
$ws.Range("F$row").Value = "NEW"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Poste podrido"
$ws.Range("I$row").Value = "1"
$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Poste"

# Coordinates stay as real numbers.
$ws.Range("M$row").Value = -58.517389
$ws.Range("N$row").Value = -34.593541
